$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6803744281114348
$ws.Range("C2").Value = 0.2145998752147662
$ws.Range("D2").Value = 0.2109788960596148
$ws.Range("F2").Value = 1.440208222474048
$ws.Range("G2").Value = 0.8003659862271277
$ws.Range("H2").Value = 0.8903937029114033
$ws.Range("J2").Value = 0.2257039721467713
$ws.Range("K2").Value = 0.3115370086762539
$ws.Range("L2").Value = 0.3221036421143424
$ws.Range("N2").Value = 1.808715983579878
$ws.Range("O2").Value = 3.399185732184108
# Row 3
$ws.Range("B3").Value = 0.6417126112047811
$ws.Range("C3").Value = 0.2152185454123341
$ws.Range("D3").Value = 0.2085745742284644
$ws.Range("F3").Value = 1.445677040594703
$ws.Range("G3").Value = 0.8044376424257678
$ws.Range("H3").Value = 0.8957523928771209
$ws.Range("J3").Value = 0.2268430369707133
$ws.Range("K3").Value = 0.2785611549095393
$ws.Range("L3").Value = 0.3168923392656069
$ws.Range("N3").Value = 1.824896324991139
$ws.Range("O3").Value = 3.418794036744046
# Row 4
$ws.Range("B4").Value = 0.6181690314968478
$ws.Range("C4").Value = 0.2156302678234034
$ws.Range("D4").Value = 0.207176592371809
$ws.Range("F4").Value = 1.449658139633421
$ws.Range("G4").Value = 0.8073435992009692
$ws.Range("H4").Value = 0.8993490439589138
$ws.Range("J4").Value = 0.2276200613426482
$ws.Range("K4").Value = 0.2583189851188905
$ws.Range("L4").Value = 0.3138163678670907
$ws.Range("N4").Value = 1.835357425023448
$ws.Range("O4").Value = 3.432327065961587
# Row 5
$ws.Range("B5").Value = 0.6086246635806276
$ws.Range("C5").Value = 0.215806084294826
$ws.Range("D5").Value = 0.2066266892620092
$ws.Range("F5").Value = 1.451437401970715
$ws.Range("G5").Value = 0.8086299432606125
$ws.Range("H5").Value = 0.9008918605251353
$ws.Range("J5").Value = 0.2279562594851114
$ws.Range("K5").Value = 0.2500720181453602
$ws.Range("L5").Value = 0.3125941770503999
$ws.Range("N5").Value = 1.839752820719053
$ws.Range("O5").Value = 3.438217781865902
# Row 6
$ws.Range("B6").Value = 0.6070428622484769
$ws.Range("C6").Value = 0.2158357646136899
$ws.Range("D6").Value = 0.2065365762645683
$ws.Range("F6").Value = 1.451742331843143
$ws.Range("G6").Value = 0.8088497108228125
$ws.Range("H6").Value = 0.9011527069784435
$ws.Range("J6").Value = 0.2280132669301516
$ws.Range("K6").Value = 0.2487027465510607
$ws.Range("L6").Value = 0.3123931280241123
$ws.Range("N6").Value = 1.840490672272521
$ws.Range("O6").Value = 3.439218644465058
# Row 7
$ws.Range("B7").Value = 0.6180401099657047
$ws.Range("C7").Value = 0.2156326063714715
$ws.Range("D7").Value = 0.2071690959387809
$ws.Range("F7").Value = 1.44968149968571
$ws.Range("G7").Value = 0.8073605336267562
$ws.Range("H7").Value = 0.8993695383950566
$ws.Range("J7").Value = 0.2276245162093105
$ws.Range("K7").Value = 0.2582077552383879
$ws.Range("L7").Value = 0.3137997580540812
$ws.Range("N7").Value = 1.835416166641009
$ws.Range("O7").Value = 3.432404987751681
# Row 8
$ws.Range("B8").Value = 0.6670038514300813
$ws.Range("C8").Value = 0.2148065962834558
$ws.Range("D8").Value = 0.2101336982366888
$ws.Range("F8").Value = 1.441964650141955
$ws.Range("G8").Value = 0.8016856780842119
$ws.Range("H8").Value = 0.8921778511572072
$ws.Range("J8").Value = 0.2260806283052368
$ws.Range("K8").Value = 0.3001662116089676
$ws.Range("L8").Value = 0.3202811889045734
$ws.Range("N8").Value = 1.814185773375616
$ws.Range("O8").Value = 3.40563692154501
# Row 9
$ws.Range("B9").Value = 0.7645375460863875
$ws.Range("C9").Value = 0.2134384318070417
$ws.Range("D9").Value = 0.2165644075683417
$ws.Range("F9").Value = 1.431768136692185
$ws.Range("G9").Value = 0.7937758734074549
$ws.Range("H9").Value = 0.8805014512505522
$ws.Range("J9").Value = 0.2236676646463387
$ws.Range("K9").Value = 0.3824654192757748
$ws.Range("L9").Value = 0.3339671333436485
$ws.Range("N9").Value = 1.776725390593253
$ws.Range("O9").Value = 3.364980176322717
# Row 10
$ws.Range("B10").Value = 0.8370858005774551
$ws.Range("C10").Value = 0.2125851328120589
$ws.Range("D10").Value = 0.2216605734043497
$ws.Range("F10").Value = 1.42727534134854
$ws.Range("G10").Value = 0.7899241092241027
$ws.Range("H10").Value = 0.873395942528191
$ws.Range("J10").Value = 0.2222677212963156
$ws.Range("K10").Value = 0.4429183948087427
$ws.Range("L10").Value = 0.3446099183828153
$ws.Range("N10").Value = 1.751741020520786
$ws.Range("O10").Value = 3.342306988292933
# Row 11
$ws.Range("B11").Value = 0.8702764262517917
$ws.Range("C11").Value = 0.2122296086751234
$ws.Range("D11").Value = 0.224058666913308
$ws.Range("F11").Value = 1.425880388135653
$ws.Range("G11").Value = 0.7885968504308636
$ws.Range("H11").Value = 0.8704820676223193
$ws.Range("J11").Value = 0.221711436039314
$ws.Range("K11").Value = 0.4704126632558712
$ws.Range("L11").Value = 0.3495777643040014
$ws.Range("N11").Value = 1.740924806738205
$ws.Range("O11").Value = 3.333551564601294
# Row 12
$ws.Range("B12").Value = 0.8828711586351972
$ws.Range("C12").Value = 0.2120996499060972
$ws.Range("D12").Value = 0.2249781462715532
$ws.Range("F12").Value = 1.425445260111559
$ws.Range("G12").Value = 0.7881553051885817
$ws.Range("H12").Value = 0.8694243506796795
$ws.Range("J12").Value = 0.2215123377509585
$ws.Range("K12").Value = 0.4808226221116172
$ws.Range("L12").Value = 0.3514769824284372
$ws.Range("N12").Value = 1.736907912927762
$ws.Range("O12").Value = 3.330459940272306
# Row 13
$ws.Range("B13").Value = 0.880157508091429
$ws.Range("C13").Value = 0.2121274315107264
$ws.Range("D13").Value = 0.2247796157776349
$ws.Range("F13").Value = 1.42553483466758
$ws.Range("G13").Value = 0.7882476849743796
$ws.Range("H13").Value = 0.869650117659944
$ws.Range("J13").Value = 0.2215547036848875
$ws.Range("K13").Value = 0.4785807293051221
$ws.Range("L13").Value = 0.3510671538906678
$ws.Range("N13").Value = 1.737769510848086
$ws.Range("O13").Value = 3.331115824878481
# Row 14
$ws.Range("B14").Value = 0.8713120831675383
$ws.Range("C14").Value = 0.212218823422063
$ws.Range("D14").Value = 0.224134085695411
$ws.Range("F14").Value = 1.425842724920471
$ws.Range("G14").Value = 0.7885593007888332
$ws.Range("H14").Value = 0.8703941331423124
$ws.Range("J14").Value = 0.2216948246924026
$ws.Range("K14").Value = 0.4712691310882349
$ws.Range("L14").Value = 0.3497336543076699
$ws.Range("N14").Value = 1.740592751848702
$ws.Range("O14").Value = 3.333292729819703
# Row 15
$ws.Range("B15").Value = 0.8658973860300989
$ws.Range("C15").Value = 0.2122754111364742
$ws.Range("D15").Value = 0.2237401579164811
$ws.Range("F15").Value = 1.426043436614371
$ws.Range("G15").Value = 0.7887581247864972
$ws.Range("H15").Value = 0.8708558135023878
$ws.Range("J15").Value = 0.2217821567802645
$ws.Range("K15").Value = 0.4667903483781117
$ws.Range("L15").Value = 0.3489191866817549
$ws.Range("N15").Value = 1.742332352443508
$ws.Range("O15").Value = 3.334655292480278
# Row 16
$ws.Range("B16").Value = 0.8349204241700647
$ws.Range("C16").Value = 0.2126090220197732
$ws.Range("D16").Value = 0.2215054496118256
$ws.Range("F16").Value = 1.427379544757912
$ws.Range("G16").Value = 0.7900193954867092
$ws.Range("H16").Value = 0.8735927716143266
$ws.Range("J16").Value = 0.2223056943596227
$ws.Range("K16").Value = 0.4411214017307259
$ws.Range("L16").Value = 0.3442877865550997
$ws.Range("N16").Value = 1.752458937347512
$ws.Range("O16").Value = 3.342910517569322
# Row 17
$ws.Range("B17").Value = 0.8159646022158427
$ws.Range("C17").Value = 0.2128220257940505
$ws.Range("D17").Value = 0.220154898176034
$ws.Range("F17").Value = 1.42836525115807
$ws.Range("G17").Value = 0.7909019464129443
$ws.Range("H17").Value = 0.8753533073975319
$ws.Range("J17").Value = 0.2226474793841255
$ws.Range("K17").Value = 0.4253723129653508
$ws.Range("L17").Value = 0.3414788217102966
$ws.Range("N17").Value = 1.758811945419117
$ws.Range("O17").Value = 3.348373859545262
# Row 18
$ws.Range("B18").Value = 0.8050794686075449
$ws.Range("C18").Value = 0.2129476141343929
$ws.Range("D18").Value = 0.2193856162272283
$ws.Range("F18").Value = 1.428993286991322
$ws.Range("G18").Value = 0.7914495653392919
$ws.Range("H18").Value = 0.8763959009842637
$ws.Range("J18").Value = 0.2228516491316341
$ws.Range("K18").Value = 0.4163133188512802
$ws.Range("L18").Value = 0.3398750861140485
$ws.Range("N18").Value = 1.762517731364422
$ws.Range("O18").Value = 3.351662960387785
# Row 19
$ws.Range("B19").Value = 0.8013970254463914
$ws.Range("C19").Value = 0.212990664917097
$ws.Range("D19").Value = 0.2191264452563075
$ws.Range("F19").Value = 1.429216426149175
$ws.Range("G19").Value = 0.7916418509334733
$ws.Range("H19").Value = 0.8767540570791255
$ws.Range("J19").Value = 0.2229220809504717
$ws.Range("K19").Value = 0.4132460269781859
$ws.Range("L19").Value = 0.3393341390183338
$ws.Range("N19").Value = 1.763781327544519
$ws.Range("O19").Value = 3.35280180308169
# Row 20
$ws.Range("B20").Value = 0.817980650743948
$ws.Range("C20").Value = 0.2127990331845169
$ws.Range("D20").Value = 0.2202978892015324
$ws.Range("F20").Value = 1.428254000498931
$ws.Range("G20").Value = 0.7908038579426631
$ws.Range("H20").Value = 0.8751627931924446
$ws.Range("J20").Value = 0.222610311075421
$ws.Range("K20").Value = 0.4270488903910348
$ws.Range("L20").Value = 0.3417766097612827
$ws.Range("N20").Value = 1.758130305854374
$ws.Range("O20").Value = 3.347777093365181
# Row 21
$ws.Range("B21").Value = 0.873909495850711
$ws.Range("C21").Value = 0.2121918528352609
$ws.Range("D21").Value = 0.2243233856944045
$ws.Range("F21").Value = 1.425749764631547
$ws.Range("G21").Value = 0.7884661149138026
$ws.Range("H21").Value = 0.8701743580706705
$ws.Range("J21").Value = 0.2216533543908064
$ws.Range("K21").Value = 0.4734167712418582
$ws.Range("L21").Value = 0.3501248482435528
$ws.Range("N21").Value = 1.739761354389813
$ws.Range("O21").Value = 3.332647246664919
# Row 22
$ws.Range("B22").Value = 0.9106143171099745
$ws.Range("C22").Value = 0.2118222376838403
$ws.Range("D22").Value = 0.227020524419828
$ws.Range("F22").Value = 1.424655731453932
$ws.Range("G22").Value = 0.7872941410902285
$ws.Range("H22").Value = 0.8671804854314615
$ws.Range("J22").Value = 0.2210952650701152
$ws.Range("K22").Value = 0.5037116785979379
$ws.Range("L22").Value = 0.355685758156838
$ws.Range("N22").Value = 1.728216483356139
$ws.Range("O22").Value = 3.324063742286967
# Row 23
$ws.Range("B23").Value = 0.8910106381121921
$ws.Range("C23").Value = 0.2120170261790761
$ws.Range("D23").Value = 0.2255749840813053
$ws.Range("F23").Value = 1.425190050481568
$ws.Range("G23").Value = 0.7878870979512413
$ws.Range("H23").Value = 0.8687540292721678
$ws.Range("J23").Value = 0.2213869760480307
$ws.Range("K23").Value = 0.4875437707256651
$ws.Range("L23").Value = 0.3527082586686276
$ws.Range("N23").Value = 1.734336085856743
$ws.Range("O23").Value = 3.328525629580128
# Row 24
$ws.Range("B24").Value = 0.8170691552017217
$ws.Range("C24").Value = 0.2128094183942757
$ws.Range("D24").Value = 0.220233220623939
$ws.Range("F24").Value = 1.428304105838905
$ws.Range("G24").Value = 0.7908480784278993
$ws.Range("H24").Value = 0.8752488298455035
$ws.Range("J24").Value = 0.2226270909654069
$ws.Range("K24").Value = 0.4262909241194563
$ws.Range("L24").Value = 0.3416419449603438
$ws.Range("N24").Value = 1.758438308896631
$ws.Range("O24").Value = 3.348046429850399
# Row 25
$ws.Range("B25").Value = 0.7379933094377407
$ws.Range("C25").Value = 0.213781774095807
$ws.Range("D25").Value = 0.2147591472327699
$ws.Range("F25").Value = 1.433999210565005
$ws.Range("G25").Value = 0.7955713782951435
$ws.Range("H25").Value = 0.8834010883666039
$ws.Range("J25").Value = 0.2242548262523307
$ws.Range("K25").Value = 0.3602017410856035
$ws.Range("L25").Value = 0.3301609918094783
$ws.Range("N25").Value = 1.786413347356612
$ws.Range("O25").Value = 3.374713690945043
